# Auto-applied crypto price/volume update, commit: 'Updated cryptos list on Fri Apr 12 19:19:38 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level Price (D) / Volume(1h) (E) updates
$updates = @(
    @{ Row=2; D="66.495.65"; E="  -5.18%  " },
    @{ Row=3; D="3.217.00"; E="  -8.17%  " },
    @{ Row=4; D="1.00"; E="  -0.02%  " },
    @{ Row=5; D="587.50"; E="  -2.63%  " },
    @{ Row=6; D="151.77"; E="  -13.36%  " },
    @{ Row=7; D="0.999"; E="  -0.14%  " },
    @{ Row=8; D="3.212.26"; E="  -8.14%  " },
    @{ Row=9; D="0.527"; E="  -13.68%  " },
    @{ Row=10; D="0.172"; E="  -10.80%  " },
    @{ Row=11; D="6.38"; E="  -11.21%  " },
    @{ Row=12; D="0.481"; E="  -17.08%  " },
    @{ Row=13; D="38.54"; E="  -16.45%  " },
    @{ Row=14; D=$null; E="  -13.51%  " },
    @{ Row=15; D="3.744.65"; E="  -7.77%  " },
    @{ Row=16; D="66.539.49"; E="  -5.24%  " },
    @{ Row=17; D="3.217.45"; E="  -8.17%  " },
    @{ Row=18; D="0.114"; E="  -4.65%  " },
    @{ Row=19; D="514.94"; E="  -15.56%  " },
    @{ Row=20; D=$null; E="  -16.74%  " },
    @{ Row=21; D="14.27"; E="  -17.55%  " },
    @{ Row=22; D="0.742"; E="  -15.08%  " },
    @{ Row=23; D="7.75"; E="  -13.67%  " },
    @{ Row=24; D="84.08"; E="  -14.45%  " },
    @{ Row=25; D="13.08"; E="  -15.71%  " },
    @{ Row=26; D="1.00"; E="  -0.02%  " },
    @{ Row=27; D="3.18"; E="  -14.38%  " },
    @{ Row=28; D=$null; E="  -17.34%  " },
    @{ Row=29; D="28.54"; E="  -15.59%  " },
    @{ Row=30; D=$null; E="  -16.20%  " },
    @{ Row=31; D=$null; E="  -10.79%  " },
    @{ Row=32; D=$null; E="  -10.62%  " },
    @{ Row=33; D="539.29"; E="  -14.47%  " },
    @{ Row=34; D=$null; E="  -16.39%  " },
    @{ Row=35; D="6.43"; E="  -19.80%  " },
    @{ Row=36; D=$null; E="  +0.38%  " },
    @{ Row=37; D="53.38"; E="  -5.78%  " },
    @{ Row=38; D="0.0422"; E="  -10.85%  " },
    @{ Row=39; D="0.0845"; E="  -14.80%  " },
    @{ Row=40; D="9.06"; E="  -15.51%  " },
    @{ Row=41; D=$null; E="  -13.34%  " },
    @{ Row=42; D=$null; E="  -21.70%  " },
    @{ Row=43; D="2.880.24"; E="  -14.30%  " },
    @{ Row=49; D="25.51"; E="  -20.56%  " },
    @{ Row=50; D="0.111"; E="  -13.87%  " },
    @{ Row=51; D="122.04"; E="  -8.02%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value2 = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value2 = $u.E
    }
}

# Rows 44/45 and 47/48 swap position (ranking reshuffled) with refreshed values
$fullRows = @(
    @{ Row=44; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="0.0₃0572"; E="  -21.67%  " },
    @{ Row=45; B="ThetaToken"; C="https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"; D="2.40"; E="  -16.80%  " },
    @{ Row=47; B="USDe"; C="https://coinranking.com/coin/exbfr2U-0+usde-usde"; D="1.00"; E="  -0.04%  " },
    @{ Row=48; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="2.14"; E="  -16.15%  " }
)

foreach ($u in $fullRows) {
    $ws.Cells.Item($u.Row, 2).Value2 = $u.B
    $ws.Cells.Item($u.Row, 3).Value2 = $u.C
    $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
    $ws.Cells.Item($u.Row, 4).Value2 = $u.D
    $ws.Cells.Item($u.Row, 5).Value2 = $u.E
}
